$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5449382066726685
$ws.Range("B1").Value = 1.926955699920654
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.338976621627808
$ws.Range("E1").Value = 1.346160173416138
